# CIV-6625 Update GA order template
# Remove the "Classification: Controlled" text-box drawing (anchored shape)
# from the default (primary) footer of the document's only section.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 -> this is the default footer (footer2.xml),
# which contains the "Text Box 4" classification marker shape.
$footer = $sec.Footers.Item(1)

for ($i = $footer.Shapes.Count; $i -ge 1; $i--) {
    $shape = $footer.Shapes.Item($i)
    if ($shape.Name -eq "Text Box 4") {
        $shape.Delete()
    }
}
